# Add a new customer row (phone 79174463) with 0 points, mirroring the
# existing "no birthday on file" rows (e.g. row 4: 79174445).
#
# Column A stores the phone number as text (not a number) for this new
# row, and column B is left as an empty text cell (matching the blank
# "birthday" cells already present for several other rows), so we force
# text formatting before writing, then drop back to the default "Normal"
# style so no stray number-format style is left on the cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 15

# A15: phone number, stored as text.
$ws.Range("A$newRow").NumberFormat = "@"
$ws.Range("A$newRow").Value = "79174463"
$ws.Range("A$newRow").Style = "Normal"

# B15: no birthday on file -> empty text cell (a bare "'" writes an
# empty-string text value instead of clearing the cell entirely).
$ws.Range("B$newRow").Value = "'"
$ws.Range("B$newRow").Style = "Normal"

# C15: total_points starts at 0.
$ws.Range("C$newRow").Value = 0
